$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the ISIN/Name/Currency/NAV per share data (columns A-D) for rows 2-3,
# and clear the numeric values in columns E-G for rows 2-3 while keeping
# their existing (Comma) number format.
$ws.Range("A2:D3").ClearContents()
$ws.Range("E2:G3").ClearContents()

# Update the current selection to match the template's blanked-out state.
$ws.Range("A2:XFD3").Select()
